$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1349.16
$ws.Range("I15").Value = 1349.16
$ws.Range("K15").Value = 4047.48
$ws.Range("M15").Value = -3878.48
$ws.Range("H39").Value = 705.2143
$ws.Range("I39").Value = 631.0833
$ws.Range("J39").Value = 1150
$ws.Range("K39").Value = 1893.2499
$ws.Range("L39").Value = 3450
$ws.Range("M39").Value = -1597.2499
$ws.Range("N39").Value = -4042
$ws.Range("H40").Value = 1301
$ws.Range("J40").Value = 1202
$ws.Range("L40").Value = 1202
$ws.Range("N40").Value = -1552
$ws.Range("H76").Value = 4319.2666
$ws.Range("I76").Value = 4550.75
$ws.Range("J76").Value = 4235.091
$ws.Range("K76").Value = 4550.75
$ws.Range("L76").Value = 4235.091
$ws.Range("M76").Value = -4235.75
$ws.Range("N76").Value = -4865.091
$ws.Range("H79").Value = 4319.2666
$ws.Range("I79").Value = 4550.75
$ws.Range("J79").Value = 4235.091
$ws.Range("K79").Value = 4550.75
$ws.Range("L79").Value = 4235.091
$ws.Range("M79").Value = -3458.75
$ws.Range("N79").Value = -6419.091
$ws.Range("H112").Value = 2086.7896
$ws.Range("I112").Value = 1100
$ws.Range("J112").Value = 2202.8823
$ws.Range("K112").Value = 3300
$ws.Range("L112").Value = 6608.646900000001
$ws.Range("M112").Value = -2192
$ws.Range("N112").Value = -8824.6469
$ws.Range("H127").Value = 1648.8695
$ws.Range("I127").Value = 802.1
$ws.Range("J127").Value = 2300.2307
$ws.Range("K127").Value = 2406.3
$ws.Range("L127").Value = 6900.6921
$ws.Range("M127").Value = 2553.7
$ws.Range("N127").Value = -16820.6921
$ws.Range("H137").Value = 3382.8096
$ws.Range("I137").Value = 3345.6365
$ws.Range("J137").Value = 3423.7
$ws.Range("K137").Value = 10036.9095
$ws.Range("L137").Value = 10271.1
$ws.Range("M137").Value = -7486.9095
$ws.Range("N137").Value = -15371.1
$ws.Range("H138").Value = 2593.663
$ws.Range("I138").Value = 1543.6666
$ws.Range("J138").Value = 2798.2078
$ws.Range("K138").Value = 4630.9998
$ws.Range("L138").Value = 8394.6234
$ws.Range("M138").Value = 509.0002000000004
$ws.Range("N138").Value = -18674.6234

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12209.965
$ws.Range("I32").Value = 9325.076999999999
$ws.Range("J32").Value = 16755.848
$ws.Range("K32").Value = 9325.076999999999
$ws.Range("L32").Value = 16755.848
$ws.Range("M32").Value = -9038.076999999999
$ws.Range("N32").Value = -17329.848
$ws.Range("H97").Value = 4633.375
$ws.Range("I97").Value = 471.61905
$ws.Range("J97").Value = 33765.668
$ws.Range("K97").Value = 471.61905
$ws.Range("L97").Value = 33765.668
$ws.Range("M97").Value = 24.38094999999998
$ws.Range("N97").Value = -34757.668
$ws.Range("H107").Value = 55000
$ws.Range("J107").Value = 55000
$ws.Range("L107").Value = 55000
$ws.Range("N107").Value = -62680
$ws.Range("H109").Value = 40000
$ws.Range("J109").Value = 40000
$ws.Range("L109").Value = 40000
$ws.Range("N109").Value = -42774
$ws.Range("H110").Value = 531
$ws.Range("I110").Value = 528.4666999999999
$ws.Range("K110").Value = 528.4666999999999
$ws.Range("M110").Value = 1516.5333

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 45458464
$ws.Range("I86").Value = 50003860
$ws.Range("K86").Value = 50003860
$ws.Range("M86").Value = -50002737
$ws.Range("H89").Value = 45458464
$ws.Range("I89").Value = 50003860
$ws.Range("K89").Value = 250019300
$ws.Range("M89").Value = -250013684
$ws.Range("H94").Value = 6945086
$ws.Range("I94").Value = 8621304
$ws.Range("J94").Value = 755.7143
$ws.Range("K94").Value = 8621304
$ws.Range("L94").Value = 755.7143
$ws.Range("M94").Value = -8620853
$ws.Range("N94").Value = -1657.7143
$ws.Range("H134").Value = 6480.55
$ws.Range("I134").Value = 895
$ws.Range("K134").Value = 2685
$ws.Range("M134").Value = -150

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 10002302
$ws.Range("I62").Value = 2402.2222
$ws.Range("K62").Value = 2402.2222
$ws.Range("M62").Value = -1778.2222
$ws.Range("H65").Value = 10002302
$ws.Range("I65").Value = 2402.2222
$ws.Range("K65").Value = 12011.111
$ws.Range("M65").Value = -8891.111000000001
$ws.Range("H86").Value = 5600539.5
$ws.Range("I86").Value = 8361458
$ws.Range("K86").Value = 8361458
$ws.Range("M86").Value = -8360335
$ws.Range("H89").Value = 5600539.5
$ws.Range("I89").Value = 8361458
$ws.Range("K89").Value = 41807290
$ws.Range("M89").Value = -41801674
$ws.Range("H107").Value = 1056.125
$ws.Range("I107").Value = 523.4545000000001
$ws.Range("K107").Value = 523.4545000000001
$ws.Range("M107").Value = 1396.5455
$ws.Range("H109").Value = 16000.333
$ws.Range("J109").Value = 16000.333
$ws.Range("L109").Value = 16000.333
$ws.Range("N109").Value = -18080.333
$ws.Range("H132").Value = 1877
$ws.Range("I132").Value = 1547.625
$ws.Range("J132").Value = 2931
$ws.Range("K132").Value = 4642.875
$ws.Range("L132").Value = 8793
$ws.Range("M132").Value = -2112.875
$ws.Range("N132").Value = -13853

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 152.25
$ws.Range("I38").Value = 37.5
$ws.Range("J38").Value = 267
$ws.Range("K38").Value = 112.5
$ws.Range("L38").Value = 801
$ws.Range("M38").Value = 234.5
$ws.Range("N38").Value = -1495
$ws.Range("H92").Value = 647.6
$ws.Range("I92").Value = 647.6
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 1942.8
$ws.Range("L92").Value = 0
$ws.Range("M92").ClearContents()
$ws.Range("N92").Value = -694.8000000000002
$ws.Range("H113").Value = 692.2778
$ws.Range("J113").Value = 756.7
$ws.Range("L113").Value = 2270.1
$ws.Range("N113").Value = -6610.1
$ws.Range("H131").Value = 31297918
$ws.Range("I131").Value = 166667170
$ws.Range("K131").Value = 500001510
$ws.Range("M131").Value = -499996470

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H42").Value = 38263.6
$ws.Range("J42").Value = 38263.6
$ws.Range("L42").Value = 38263.6
$ws.Range("N42").Value = -39233.6
$ws.Range("H70").Value = 23687988
$ws.Range("I70").Value = 19234546
$ws.Range("K70").Value = 19234546
$ws.Range("M70").Value = -19234276
$ws.Range("H73").Value = 23687988
$ws.Range("I73").Value = 19234546
$ws.Range("K73").Value = 19234546
$ws.Range("M73").Value = -19233610
$ws.Range("H80").Value = 3743.6365
$ws.Range("I80").Value = 1800
$ws.Range("J80").Value = 4854.2856
$ws.Range("K80").Value = 1800
$ws.Range("L80").Value = 4854.2856
$ws.Range("M80").Value = -802
$ws.Range("N80").Value = -6850.2856
$ws.Range("H83").Value = 3743.6365
$ws.Range("I83").Value = 1800
$ws.Range("J83").Value = 4854.2856
$ws.Range("K83").Value = 9000
$ws.Range("L83").Value = 24271.428
$ws.Range("M83").Value = -4008
$ws.Range("N83").Value = -34255.428
$ws.Range("H93").Value = 30000
$ws.Range("J93").Value = 30000
$ws.Range("L93").Value = 30000
$ws.Range("N93").Value = -33744
$ws.Range("H115").Value = 38263.6
$ws.Range("J115").Value = 38263.6
$ws.Range("L115").Value = 38263.6
$ws.Range("N115").Value = -40613.6

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2213.111
$ws.Range("I40").Value = 2069.1333
$ws.Range("K40").Value = 2069.1333
$ws.Range("M40").Value = -1933.1333
$ws.Range("H46").Value = 7665.8335
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").ClearContents()
$ws.Range("H93").Value = 1477
$ws.Range("I93").Value = 1267.6666
$ws.Range("J93").Value = 2105
$ws.Range("K93").Value = 1267.6666
$ws.Range("L93").Value = 2105
$ws.Range("M93").Value = -19.66660000000002
$ws.Range("N93").Value = -4601

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 848.75
$ws.Range("I96").Value = 603
$ws.Range("J96").Value = 930.6667
$ws.Range("K96").Value = 603
$ws.Range("L96").Value = 930.6667
$ws.Range("M96").Value = 770
$ws.Range("N96").Value = -3676.6667
$ws.Range("H107").Value = 327.8421
$ws.Range("I107").Value = 259.93332
$ws.Range("J107").Value = 582.5
$ws.Range("K107").Value = 779.7999599999999
$ws.Range("L107").Value = 1747.5
$ws.Range("M107").Value = 1140.20004
$ws.Range("N107").Value = -5587.5
$ws.Range("H115").Value = 35998.6
$ws.Range("I115").Value = 10000
$ws.Range("J115").Value = 42498.25
$ws.Range("K115").Value = 10000
$ws.Range("L115").Value = 42498.25
$ws.Range("M115").Value = -8433
$ws.Range("N115").Value = -45632.25
$ws.Range("H126").Value = 47620064
$ws.Range("I126").Value = 62501076
$ws.Range("J126").Value = 828
$ws.Range("K126").Value = 187503228
$ws.Range("L126").Value = 2484
$ws.Range("M126").Value = -187500758
$ws.Range("N126").Value = -7424
$ws.Range("H136").Value = 1326.1482
$ws.Range("I136").Value = 501.82352
$ws.Range("J136").Value = 2727.5
$ws.Range("K136").Value = 1505.47056
$ws.Range("L136").Value = 8182.5
$ws.Range("M136").Value = 1044.52944
$ws.Range("N136").Value = -13282.5
